$wb = $excel.ActiveWorkbook

# --- Sheet: Leveraged Free Cash Flow ---
$ws1 = $wb.Worksheets.Item("Leveraged Free Cash Flow")
$ws1.Range("B3").Value = 3287000
$ws1.Range("C3").Value = 2825000
$ws1.Range("D3").Value = 2611000
$ws1.Range("E3").Value = 2334000

$ws1.Range("B25").Value = 111000
$ws1.Range("C25").Value = -2000
$ws1.Range("D25").Value = 95000
$ws1.Range("E25").Value = 45000

$ws1.Range("B57").Value = -27214
$ws1.Range("C57").Value = 2311
$ws1.Range("D57").Value = 365763
$ws1.Range("E57").Value = -21361

$ws1.Range("B84").Value = 3005000
$ws1.Range("C84").Value = 3782000
$ws1.Range("D84").Value = 6328000
$ws1.Range("E84").Value = -1371000

$ws1.Range("B96").Value = 66894
$ws1.Range("C96").Value = 124005
$ws1.Range("D96").Value = 38251
$ws1.Range("E96").Value = 159768

# --- Sheet: Cash & Equivalents ---
$ws2 = $wb.Worksheets.Item("Cash & Equivalents")
$ws2.Range("B3").Value = 52074000
$ws2.Range("C3").Value = 50469000
$ws2.Range("D3").Value = 47065000
$ws2.Range("E3").Value = 48201000

$ws2.Range("B25").Value = 358000
$ws2.Range("C25").Value = 183000
$ws2.Range("D25").Value = 430000
$ws2.Range("E25").Value = 814000

$ws2.Range("B57").Value = 310740
$ws2.Range("C57").Value = 357122
$ws2.Range("D57").Value = 600116
$ws2.Range("E57").Value = 257525

$ws2.Range("B84").Value = 6184000
$ws2.Range("C84").Value = 6012000
$ws2.Range("D84").Value = 7281000
$ws2.Range("E84").Value = 8736000

$ws2.Range("B96").Value = 308042
$ws2.Range("C96").Value = 385790
$ws2.Range("D96").Value = 290899
$ws2.Range("E96").Value = 405709

# --- Sheet: Debt ---
$ws3 = $wb.Worksheets.Item("Debt")
$ws3.Range("B3").Value = 20909000
$ws3.Range("C3").Value = 20892000
$ws3.Range("D3").Value = 17147000
$ws3.Range("E3").Value = 16226000

$ws3.Range("B25").Value = 1332000
$ws3.Range("C25").Value = 1333000
$ws3.Range("D25").Value = 1333000
$ws3.Range("E25").Value = 1481000

$ws3.Range("B57").Value = 247437
$ws3.Range("C57").Value = 192941
$ws3.Range("D57").Value = 195726
$ws3.Range("E57").Value = 200373

$ws3.Range("B84").Value = 54175000
$ws3.Range("C84").Value = 55545000
$ws3.Range("D84").Value = 56641000
$ws3.Range("E84").Value = 61822000

$ws3.Range("B96").Value = 1269495
$ws3.Range("C96").Value = 1289661
$ws3.Range("D96").Value = 1302230
$ws3.Range("E96").Value = 1315596

# --- Active sheet / selected tab changes ---
# Before: activeTab=2 (Debt sheet selected). After: activeTab=0 (Leveraged Free Cash Flow selected)
$ws1.Activate()
